$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header updates
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 value updates
$ws.Range("B2").Value = 235.01390878421498
$ws.Range("C2").Value = 172.57869458606558
$ws.Range("D2").Value = 235.24839503447086
$ws.Range("E2").Value = 174.72133577521726

# Row 3 value updates
$ws.Range("B3").Value = 216.4499646249308
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 210.74619960290974
$ws.Range("E3").Value = 179.08676578159518

# Update the selection to reflect the saved range
$ws.Range("B1:E3").Select()
